# Change columns language in table and script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the header row from English to Russian
$ws.Range("A1").Value = "Категория"
$ws.Range("B1").Value = "Имя"
$ws.Range("C1").Value = "Сорт"
$ws.Range("D1").Value = "Цена"
$ws.Range("E1").Value = "Картинка"
$ws.Range("F1").Value = "Акция"

# Update the view: scroll/select to match the new selection state
$ws.Range("F1").Select()
$excel.ActiveWindow.ScrollColumn = 3

# Minor width tweaks
$ws.Columns.Item(3).ColumnWidth = 23.42
$ws.Cells.DefaultColumnWidth = 14.5
